$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.675.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.124.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.35%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.116.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("E10").Value = "  -2.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.01%  "

$ws.Range("E13").Value = "  -2.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.637.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.120"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.738.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.121.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "

$ws.Range("E19").Value = "  -1.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.704"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("E23").Value = "  -5.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.73%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.92%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.66%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.20%  "

$ws.Range("E34").Value = "  -2.27%  "

$ws.Range("E35").Value = "  -3.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0738"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "436.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.94%  "

$ws.Range("E41").Value = "  -1.92%  "

$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.865.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("E45").Value = "  -3.96%  "

$ws.Range("E46").Value = "  -5.76%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.13%  "

$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.35%  "
